# Weekly update: insert 3 new price records (week of 2023-11-28, serial 45258)
# at the top of the "Vega Monumental Concepción - Frutilla" data block (row 498),
# pushing the existing rows 498:594 down to 501:597.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows above the current row 498 (shifts 498:594 -> 501:597).
$ws.Rows("498:500").Insert()

# New row 498: Especial
$ws.Range("A498").Value = 11
$ws.Range("B498").Value = "Vega Monumental Concepción"
$ws.Range("C498").Value = "Bíobío"
$ws.Range("D498").Value = 45258
$ws.Range("E498").Value = 8
$ws.Range("F498").Value = "Fruta"
$ws.Range("G498").Value = 100101
$ws.Range("H498").Value = "Berries"
$ws.Range("I498").Value = 100112025
$ws.Range("J498").Value = "Frutilla"
$ws.Range("K498").Value = "Sin especificar"
$ws.Range("L498").Value = "Especial"
$ws.Range("M498").Value = 100
$ws.Range("N498").Value = 12000
$ws.Range("O498").Value = 12000
$ws.Range("P498").Value = 12000
$ws.Range("Q498").Value = "$/bandeja 7 kilos"
$ws.Range("R498").Value = "Provincia de Melipilla"
$ws.Range("S498").Value = 1714
$ws.Range("T498").Value = 7

# New row 499: Primera
$ws.Range("A499").Value = 11
$ws.Range("B499").Value = "Vega Monumental Concepción"
$ws.Range("C499").Value = "Bíobío"
$ws.Range("D499").Value = 45258
$ws.Range("E499").Value = 8
$ws.Range("F499").Value = "Fruta"
$ws.Range("G499").Value = 100101
$ws.Range("H499").Value = "Berries"
$ws.Range("I499").Value = 100112025
$ws.Range("J499").Value = "Frutilla"
$ws.Range("K499").Value = "Sin especificar"
$ws.Range("L499").Value = "Primera"
$ws.Range("M499").Value = 100
$ws.Range("N499").Value = 9000
$ws.Range("O499").Value = 9000
$ws.Range("P499").Value = 9000
$ws.Range("Q499").Value = "$/bandeja 7 kilos"
$ws.Range("R499").Value = "Provincia de Melipilla"
$ws.Range("S499").Value = 1286
$ws.Range("T499").Value = 7

# New row 500: Segunda
$ws.Range("A500").Value = 11
$ws.Range("B500").Value = "Vega Monumental Concepción"
$ws.Range("C500").Value = "Bíobío"
$ws.Range("D500").Value = 45258
$ws.Range("E500").Value = 8
$ws.Range("F500").Value = "Fruta"
$ws.Range("G500").Value = 100101
$ws.Range("H500").Value = "Berries"
$ws.Range("I500").Value = 100112025
$ws.Range("J500").Value = "Frutilla"
$ws.Range("K500").Value = "Sin especificar"
$ws.Range("L500").Value = "Segunda"
$ws.Range("M500").Value = 50
$ws.Range("N500").Value = 7000
$ws.Range("O500").Value = 7000
$ws.Range("P500").Value = 7000
$ws.Range("Q500").Value = "$/bandeja 7 kilos"
$ws.Range("R500").Value = "Provincia de Melipilla"
$ws.Range("S500").Value = 1000
$ws.Range("T500").Value = 7
